# Add the new "2019_monthly" worksheet with termination-by-month data,
# positioned right before the "2019_procedure_location_wrksht" sheet
# (pushing it and "2019_gestation_weeks" one slot later, matching the
# sheetId/rId renumbering seen in the target workbook).

$wb = $excel.ActiveWorkbook

$target = $wb.Worksheets.Item("2019_procedure_location_wrksht")
$monthly = $wb.Worksheets.Add($target)
$monthly.Name = "2019_monthly"

$monthlyData = @(
    @("Month",     "IN Resident",  "Nonâ€Resident"),
    @("January",   "675 (8.84%)",  "61 (.80%)"),
    @("February",  "622 (8.14%)",  "72 (.94%)"),
    @("March",     "646 (8.46%)",  "86 (1.13%)"),
    @("April",     "578 (7.57%)",  "71 (.93%)"),
    @("May",       "668 (8.75%)",  "71 (.93%)"),
    @("June",      "504 (6.60%)",  "42 (.55%)"),
    @("July",      "561 (7.35%)",  "39 (.51%)"),
    @("August",    "634 (8.30%)",  "39 (.51%)"),
    @("September", "468 (6.13%)",  "35 (.46%)"),
    @("October",   "556 (7.28%)",  "43 (.56%)"),
    @("November",  "596 (7.80%)",  "33 (.43%)"),
    @("December",  "511 (6.69%)",  "26 (.34%)")
)

for ($i = 0; $i -lt $monthlyData.Length; $i++) {
    $row = $i + 1
    $monthly.Range("A$row").Value = $monthlyData[$i][0]
    $monthly.Range("B$row").Value = $monthlyData[$i][1]
    $monthly.Range("C$row").Value = $monthlyData[$i][2]
}

$monthly.Activate()
